$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "SA6"
$ws.Range("B18").Value = "October 09, 2024"

$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat

$ws.Range("B19").Select()
